$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54, shifting existing rows 54..175 down to 55..176
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new weekly record.
# Static columns mirror the values of the (now shifted) neighboring row.
$ws.Cells.Item(54, 1).Value = 9
$ws.Cells.Item(54, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(54, 3).Value = "Metropolitana"
$ws.Cells.Item(54, 4).Value = 44498
$ws.Cells.Item(54, 5).Value = 13
$ws.Cells.Item(54, 6).Value = 300000001
$ws.Cells.Item(54, 7).Value = "Rabanito"
$ws.Cells.Item(54, 8).Value = "Sin especificar"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 7700
$ws.Cells.Item(54, 11).Value = 3500
$ws.Cells.Item(54, 12).Value = 4000
$ws.Cells.Item(54, 13).Value = 3773
$ws.Cells.Item(54, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(54, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(54, 16).Value = 38
$ws.Cells.Item(54, 17).Value = 100
$ws.Cells.Item(54, 18).Value = "Hortaliza"
